# Applies updated crypto price / 1h-volume figures to the cryptos sheet,
# per the Sat Mar 30 02:51:23 UTC 2024 GitHub Actions data refresh.
# (ThetaToken/Stellar also swap rank positions 45/46 in this refresh.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.828.33"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3
$ws.Range("D3").Value = "3.504.57"
$ws.Range("E3").Value = "  -1.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "606.42"
$c.ClearFormats()
$ws.Range("E5").Value = "  +3.20%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "191.68"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.10%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.627"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.87%  "

# Row 9
$ws.Range("E9").Value = "  -0.72%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.663"
$c.ClearFormats()
$ws.Range("E10").Value = "  +2.97%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "53.44"
$c.ClearFormats()
$ws.Range("E11").Value = "  -1.16%  "

# Row 12
$ws.Range("E12").Value = "  -0.80%  "

# Row 13
$ws.Range("E13").Value = "  +2.57%  "

# Row 14
$ws.Range("D14").Value = "4.065.32"
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "617.97"
$c.ClearFormats()
$ws.Range("E15").Value = "  +10.39%  "

# Row 16
$ws.Range("D16").Value = "69.943.01"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "12.70"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.25%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "18.87"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.19%  "

# Row 19
$ws.Range("D19").Value = "3.506.54"
$ws.Range("E19").Value = "  -2.39%  "

# Row 20
$ws.Range("E20").Value = "  -0.14%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.990"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "17.94"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "105.99"
$c.ClearFormats()
$ws.Range("E23").Value = "  +13.16%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.64"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.58%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "4.99"
$c.ClearFormats()
$ws.Range("E25").Value = "  +2.56%  "

# Row 26
$ws.Range("E26").Value = "  +4.19%  "

# Row 27
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.89"
$c.ClearFormats()
$ws.Range("E28").Value = "  +6.25%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "34.10"
$c.ClearFormats()
$ws.Range("E29").Value = "  +5.64%  "

# Row 30
$ws.Range("E30").Value = "  +0.99%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "12.67"
$c.ClearFormats()
$ws.Range("E31").Value = "  +4.16%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.17"
$c.ClearFormats()
$ws.Range("E32").Value = "  +4.83%  "

# Row 33
$ws.Range("E33").Value = "  +0.26%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "64.14"
$c.ClearFormats()
$ws.Range("E34").Value = "  +1.59%  "

# Row 35
$ws.Range("D35").Value = "3.722.54"
$ws.Range("E35").Value = "  +1.62%  "

# Row 36
$ws.Range("E36").Value = "  -4.06%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "525.18"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.95%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0797"
$ws.Range("E39").Value = "  +1.56%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.390"
$c.ClearFormats()
$ws.Range("E40").Value = "  -3.50%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.78"
$c.ClearFormats()
$ws.Range("E41").Value = "  -3.29%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.58"
$c.ClearFormats()
$ws.Range("E42").Value = "  +0.76%  "

# Row 43
$ws.Range("E43").Value = "  +0.37%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0462"
$c.ClearFormats()
$ws.Range("E44").Value = "  +1.21%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.141"
$c.ClearFormats()
$ws.Range("E45").Value = "  +2.66%  "

# Row 46
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.86"
$c.ClearFormats()
$ws.Range("E46").Value = "  -2.53%  "

# Row 47
$ws.Range("E47").Value = "  -4.34%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.76"
$c.ClearFormats()
$ws.Range("E48").Value = "  -4.46%  "

# Row 49
$ws.Range("E49").Value = "  +0.43%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "132.53"
$c.ClearFormats()
$ws.Range("E50").Value = "  -1.31%  "

# Row 51
$ws.Range("E51").Value = "  -6.80%  "

